$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-adjusted timestamp value on the existing last row (row 12)
$ws.Cells.Item(12, 1).Value = 44325.77076643982

# Append the new data row (row 13)
$ws.Cells.Item(13, 1).Value = 44326.78153125344
$ws.Cells.Item(13, 2).Value = 74047
$ws.Cells.Item(13, 3).Value = 62273
$ws.Cells.Item(13, 4).Value = 3245
$ws.Cells.Item(13, 5).Value = 2048
$ws.Cells.Item(13, 6).Value = 1439
$ws.Cells.Item(13, 7).Value = 19198
$ws.Cells.Item(13, 8).Value = 1359
$ws.Cells.Item(13, 9).Value = 832
$ws.Cells.Item(13, 10).Value = 215
